$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look numeric,
# so they are stored as text (matching the workbook's inlineStr cells)
# instead of being parsed into numbers.
$textCells = @('D4', 'D5', 'D6', 'D11', 'D14', 'D17', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D27', 'D30', 'D33', 'D34', 'D35', 'D36', 'D38', 'D39', 'D40', 'D42', 'D43', 'D44', 'D46', 'D47', 'D48', 'D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range('D2').Value = '57.034.67'
$ws.Range('E2').Value = '  -1.59%  '
$ws.Range('D3').Value = '2.984.80'
$ws.Range('E3').Value = '  -2.17%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '501.55'
$ws.Range('E5').Value = '  -4.44%  '
$ws.Range('D6').Value = '138.38'
$ws.Range('E6').Value = '  -3.05%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -3.00%  '
$ws.Range('E9').Value = '  -4.02%  '
$ws.Range('E10').Value = '  -2.83%  '
$ws.Range('D11').Value = '0.360'
$ws.Range('E11').Value = '  -2.03%  '
$ws.Range('D12').Value = '3.487.07'
$ws.Range('E12').Value = '  -2.09%  '
$ws.Range('E13').Value = '  -2.37%  '
$ws.Range('D14').Value = '26.09'
$ws.Range('E14').Value = '  -1.14%  '
$ws.Range('E15').Value = '  -3.12%  '
$ws.Range('D16').Value = '57.098.03'
$ws.Range('E16').Value = '  -1.27%  '
$ws.Range('D17').Value = '6.06'
$ws.Range('E17').Value = '  -1.15%  '
$ws.Range('D18').Value = '2.986.15'
$ws.Range('E18').Value = '  -1.78%  '
$ws.Range('D19').Value = '12.66'
$ws.Range('E19').Value = '  -2.51%  '
$ws.Range('D20').Value = '7.89'
$ws.Range('E20').Value = '  -3.14%  '
$ws.Range('D21').Value = '322.50'
$ws.Range('E21').Value = '  -5.32%  '
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').Value = '5.76'
$ws.Range('E23').Value = '  +0.23%  '
$ws.Range('D24').Value = '0.490'
$ws.Range('E24').Value = '  -1.59%  '
$ws.Range('D25').Value = '63.82'
$ws.Range('E25').Value = '  -2.02%  '
$ws.Range('E26').Value = '  +0.43%  '
$ws.Range('D27').Value = '0.164'
$ws.Range('E27').Value = '  -5.83%  '
$ws.Range('D28').Value = '0.0₃0899'
$ws.Range('E28').Value = '  -6.08%  '
$ws.Range('E29').Value = '  -5.89%  '
$ws.Range('D30').Value = '7.08'
$ws.Range('E30').Value = '  -2.50%  '
$ws.Range('E31').Value = '  -4.14%  '
$ws.Range('E32').Value = '  -5.74%  '
$ws.Range('D33').Value = '20.25'
$ws.Range('E33').Value = '  -3.66%  '
$ws.Range('D34').Value = '155.37'
$ws.Range('E34').Value = '  -2.02%  '
$ws.Range('D35').Value = '4.59'
$ws.Range('E35').Value = '  -2.95%  '
$ws.Range('D36').Value = '5.80'
$ws.Range('E36').Value = '  -1.73%  '
$ws.Range('E37').Value = '  -5.30%  '
$ws.Range('D38').Value = '24.09'
$ws.Range('E38').Value = '  -5.31%  '
$ws.Range('D39').Value = '0.0668'
$ws.Range('E39').Value = '  -3.74%  '
$ws.Range('D40').Value = '37.90'
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('D41').Value = '3.013.97'
$ws.Range('E41').Value = '  -2.10%  '
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.15%  '
$ws.Range('D43').Value = '3.75'
$ws.Range('E43').Value = '  -2.16%  '
$ws.Range('D44').Value = '0.641'
$ws.Range('E44').Value = '  -2.96%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.202.13'
$ws.Range('E45').Value = '  -6.01%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').Value = '1.39'
$ws.Range('E46').Value = '  -5.03%  '
$ws.Range('D47').Value = '0.949'
$ws.Range('E47').Value = '  -8.19%  '
$ws.Range('D48').Value = '5.98'
$ws.Range('E48').Value = '  -0.35%  '
$ws.Range('E49').Value = '  -4.54%  '
$ws.Range('D50').Value = '19.29'
$ws.Range('E50').Value = '  -3.61%  '
$ws.Range('E51').Value = '  -10.44%  '
